$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Write a string that would otherwise be auto-coerced to a number,
    # while keeping the cell General-formatted afterwards (matches source file).
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "28.326.02"
$ws.Range("E2").Value = "  +4.04%  "
$ws.Range("D3").Value = "1.581.85"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.88%  "
Set-TextValue $ws.Range("D5") "213.67"
$ws.Range("E5").Value = "  +0.97%  "
Set-TextValue $ws.Range("D6") "0.496"
$ws.Range("E6").Value = "  +0.50%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.88%  "
Set-TextValue $ws.Range("D8") "23.62"
$ws.Range("E8").Value = "  +7.02%  "
$ws.Range("E9").Value = "  +1.02%  "
Set-TextValue $ws.Range("D10") "0.0599"
$ws.Range("E10").Value = "  -0.21%  "
Set-TextValue $ws.Range("D11") "0.0886"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "1.806.91"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "1.596.53"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D16").Value = "28.300.02"
$ws.Range("E16").Value = "  +4.08%  "
Set-TextValue $ws.Range("D17") "63.84"
$ws.Range("E17").Value = "  +2.46%  "
Set-TextValue $ws.Range("D18") "232.23"
$ws.Range("E18").Value = "  +7.36%  "
$ws.Range("D19").Value = "0.0₃0709"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  -0.93%  "
Set-TextValue $ws.Range("D22") "4.13"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  +1.14%  "
Set-TextValue $ws.Range("D24") "1.94"
$ws.Range("E24").Value = "  -0.32%  "
Set-TextValue $ws.Range("D25") "151.58"
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("E28").Value = "  -0.24%  "
Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  -0.86%  "
Set-TextValue $ws.Range("D30") "1.15"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("E32").Value = "  -0.38%  "
Set-TextValue $ws.Range("D33") "3.15"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "1.418.87"
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("E38").Value = "  -0.33%  "
Set-TextValue $ws.Range("D39") "2.51"
$ws.Range("E39").Value = "  +7.21%  "
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("E41").Value = "  +0.04%  "
Set-TextValue $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("E44").Value = "  +5.41%  "
Set-TextValue $ws.Range("D45") "0.972"
$ws.Range("E45").Value = "  -2.73%  "
Set-TextValue $ws.Range("D46") "64.26"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "1.717.45"
$ws.Range("E47").Value = "  +0.81%  "
Set-TextValue $ws.Range("D48") "87.23"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("E50").Value = "  +0.44%  "
Set-TextValue $ws.Range("D51") "39.23"
$ws.Range("E51").Value = "  +15.74%  "
